# Adds per-year data for 2012-2017 (columns K:P) to the existing
# 2018-2024 table (columns D:J), matching a fresh upload of
# "data_2012-2024" sourced figures. Also:
#  - J1 changes from the shared "2024" label to a literal numeric year
#    (loses the bold/bordered header style, matching the K1:P1 cells)
#  - E15 (previously blank) gets a single-space placeholder
#  - D25:H25 (previously blank) get "-" placeholders
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header ---
# J1 used to be a styled shared-string "2024"; it becomes a plain literal number.
$ws.Range("J1").ClearFormats()
$ws.Range("J1").Value = 2024

# New year columns K:P (2012-2017), unstyled like the new J1.
$ws.Range("K1").Value = 2012
$ws.Range("L1").Value = 2013
$ws.Range("M1").Value = 2014
$ws.Range("N1").Value = 2015
$ws.Range("O1").Value = 2016
$ws.Range("P1").Value = 2017

# --- Data rows 2-35: columns K:P (years 2012-2017) ---
# Cells with no reported figure are written as "-", mirroring the source data.
# Row 2: 营业总收入
$ws.Range("K2").Value = 1908.894871794872
$ws.Range("L2").Value = 2155.5867256637166
$ws.Range("M2").Value = 2236.3767857142857
$ws.Range("N2").Value = 2077.4663636363634
$ws.Range("O2").Value = 2222.381904761905
$ws.Range("P2").Value = 2669.4363636363637
# Row 3: 营业总成本
$ws.Range("K3").Value = 1830.6145299145298
$ws.Range("L3").Value = 2069.5070796460177
$ws.Range("M3").Value = 2151.3535714285713
$ws.Range("N3").Value = 2018.6572727272726
$ws.Range("O3").Value = 2154.0828571428574
$ws.Range("P3").Value = 2573.880808080808
# Row 4: 营业利润
$ws.Range("K4").Value = 96.35726495726495
$ws.Range("L4").Value = 103.1283185840708
$ws.Range("M4").Value = 105.21160714285715
$ws.Range("N4").Value = 90.9609090909091
$ws.Range("O4").Value = 96.6647619047619
$ws.Range("P4").Value = 145.34848484848484
# Row 5: 利润总额
$ws.Range("K5").Value = 108.26581196581196
$ws.Range("L5").Value = 116.17610619469026
$ws.Range("M5").Value = 120.31428571428572
$ws.Range("N5").Value = 111.51272727272728
$ws.Range("O5").Value = 117.68380952380951
$ws.Range("P5").Value = 145.56363636363636
# Row 6: 净利润
$ws.Range("K6").Value = 78.66153846153846
$ws.Range("L6").Value = 83.03362831858406
$ws.Range("M6").Value = 86.64642857142857
$ws.Range("N6").Value = 81.13272727272728
$ws.Range("O6").Value = 83.55619047619047
$ws.Range("P6").Value = 104.35858585858585
# Row 7: 企业年末从业人数
$ws.Range("K7").Value = 10.841880341880342
$ws.Range("L7").Value = 11.393805309734514
$ws.Range("M7").Value = 12.999107142857143
$ws.Range("N7").Value = 12.896363636363635
$ws.Range("O7").Value = 12.939047619047619
$ws.Range("P7").Value = 13.56060606060606
# Row 8: 科技人员人数
$ws.Range("K8").Value = "-"
$ws.Range("L8").Value = "-"
$ws.Range("M8").Value = "-"
$ws.Range("N8").Value = "-"
$ws.Range("O8").Value = "-"
$ws.Range("P8").Value = "-"
# Row 9: 党员人数
$ws.Range("K9").Value = "-"
$ws.Range("L9").Value = "-"
$ws.Range("M9").Value = "-"
$ws.Range("N9").Value = "-"
$ws.Range("O9").Value = "-"
$ws.Range("P9").Value = 3.9808080808080812
# Row 10: 实际发放职工薪酬总额
$ws.Range("K10").Value = 82.64957264957265
$ws.Range("L10").Value = 94.97345132743362
$ws.Range("M10").Value = 166.58660714285716
$ws.Range("N10").Value = 177.77363636363634
$ws.Range("O10").Value = 194.53142857142856
$ws.Range("P10").Value = 223.1010101010101
# Row 11: 增加值（劳动生产总值）
$ws.Range("K11").Value = 416.9940170940171
$ws.Range("L11").Value = "-"
$ws.Range("M11").Value = "-"
$ws.Range("N11").Value = 529.73
$ws.Range("O11").Value = 584.4209523809524
$ws.Range("P11").Value = 694.8626262626262
# Row 12: 本年度研发经费投入
$ws.Range("K12").Value = 32.78034188034188
$ws.Range("L12").Value = 37.901769911504424
$ws.Range("M12").Value = 40.973214285714285
$ws.Range("N12").Value = 44.87
$ws.Range("O12").Value = 50.08761904761904
$ws.Range("P12").Value = 61.41414141414141
# Row 13: 累计拥有专利数
$ws.Range("K13").Value = 0.15920256410256411
$ws.Range("L13").Value = 0.2257637168141593
$ws.Range("M13").Value = 0.3035267857142857
$ws.Range("N13").Value = 0.37272727272727274
$ws.Range("O13").Value = 0.4647619047619047
$ws.Range("P13").Value = 0.5808080808080808
# Row 14: 实际上缴税费总额
$ws.Range("K14").Value = "-"
$ws.Range("L14").Value = "-"
$ws.Range("M14").Value = "-"
$ws.Range("N14").Value = 179.50272727272727
$ws.Range("O14").Value = 184.46190476190475
$ws.Range("P14").Value = 199.27777777777777
# Row 15: 对外捐赠支出总额
$ws.Range("K15").Value = "-"
$ws.Range("L15").Value = "-"
$ws.Range("M15").Value = "-"
$ws.Range("N15").Value = 0.3009090909090909
$ws.Range("O15").Value = 0.5523809523809524
$ws.Range("P15").Value = 0.6929292929292928
# Row 16: 净资产收益率
$ws.Range("K16").Value = 8.2
$ws.Range("L16").Value = 7.6
$ws.Range("M16").Value = 7.2
$ws.Range("N16").Value = 6
$ws.Range("O16").Value = 5.4
$ws.Range("P16").Value = 5.9
# Row 17: 国有资本回报率
$ws.Range("K17").Value = "-"
$ws.Range("L17").Value = "-"
$ws.Range("M17").Value = "-"
$ws.Range("N17").Value = "-"
$ws.Range("O17").Value = "-"
$ws.Range("P17").Value = "-"
# Row 18: 总资产报酬率
$ws.Range("K18").Value = 5.7
$ws.Range("L18").Value = 5.3
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = 4
$ws.Range("O18").Value = 3.6
$ws.Range("P18").Value = 3.8
# Row 19: 营业收入利润率
$ws.Range("K19").Value = 5
$ws.Range("L19").Value = 4.8
$ws.Range("M19").Value = 4.7
$ws.Range("N19").Value = 4.4
$ws.Range("O19").Value = 4.3
$ws.Range("P19").Value = 5.4
# Row 20: 成本费用利润率
$ws.Range("K20").Value = 5.9
$ws.Range("L20").Value = 5.6
$ws.Range("M20").Value = 5.6
$ws.Range("N20").Value = 5.6
$ws.Range("O20").Value = 5.5
$ws.Range("P20").Value = 5.8
# Row 21: 盈余现金保障倍数
$ws.Range("K21").Value = 2.1
$ws.Range("L21").Value = 2.2
$ws.Range("M21").Value = 2.4
$ws.Range("N21").Value = 2.9
$ws.Range("O21").Value = 3
$ws.Range("P21").Value = 2.6
# Row 22: 成本费用总额占营业总收入的比率
$ws.Range("K22").Value = 95.4
$ws.Range("L22").Value = 95.4
$ws.Range("M22").Value = 95.4
$ws.Range("N22").Value = 96
$ws.Range("O22").Value = 95.5
$ws.Range("P22").Value = 94.8
# Row 23: 总资产增长率
$ws.Range("K23").Value = 11.7
$ws.Range("L23").Value = 11.2
$ws.Range("M23").Value = 10.1
$ws.Range("N23").Value = 22.3
$ws.Range("O23").Value = 7.4
$ws.Range("P23").Value = 8
# Row 24: 应收账款增长率
$ws.Range("K24").Value = 19.5
$ws.Range("L24").Value = 12.9
$ws.Range("M24").Value = 12.6
$ws.Range("N24").Value = 2.5
$ws.Range("O24").Value = 8.4
$ws.Range("P24").Value = 8.3
# Row 25: 营业现金比率
$ws.Range("K25").Value = "-"
$ws.Range("L25").Value = "-"
$ws.Range("M25").Value = "-"
$ws.Range("N25").Value = "-"
$ws.Range("O25").Value = "-"
$ws.Range("P25").Value = "-"
# Row 26: 资产负债率
$ws.Range("K26").Value = 62.7
$ws.Range("L26").Value = 63.4
$ws.Range("M26").Value = 63
$ws.Range("N26").Value = 66.7
$ws.Range("O26").Value = 66.7
$ws.Range("P26").Value = 66.2
# Row 27: 已获利息倍数
$ws.Range("K27").Value = 4
$ws.Range("L27").Value = 3.9
$ws.Range("M27").Value = 3.7
$ws.Range("N27").Value = 3.4
$ws.Range("O27").Value = 3.5
$ws.Range("P27").Value = 3.7
# Row 28: 资本积累率
$ws.Range("K28").Value = 9.8
$ws.Range("L28").Value = 9.5
$ws.Range("M28").Value = 11.4
$ws.Range("N28").Value = 12.8
$ws.Range("O28").Value = 7
$ws.Range("P28").Value = 9.4
# Row 29: 研发经费投入强度
$ws.Range("K29").Value = 1.7
$ws.Range("L29").Value = 1.8
$ws.Range("M29").Value = 1.8
$ws.Range("N29").Value = 2.2
$ws.Range("O29").Value = 2.3
$ws.Range("P29").Value = 2.3
# Row 30: 人均资产
$ws.Range("K30").Value = 247.2
$ws.Range("L30").Value = 272
$ws.Range("M30").Value = 265.6
$ws.Range("N30").Value = 335.4
$ws.Range("O30").Value = 371.6
$ws.Range("P30").Value = 406.6
# Row 31: 人均利润
$ws.Range("K31").Value = 9.9
$ws.Range("L31").Value = 10.2
$ws.Range("M31").Value = 9.2
$ws.Range("N31").Value = 8.5
$ws.Range("O31").Value = 9
$ws.Range("P31").Value = 10.6
# Row 32: 职工人均工资
$ws.Range("K32").Value = 7.6
$ws.Range("L32").Value = 8.3
$ws.Range("M32").Value = 8.9
$ws.Range("N32").Value = 9.4
$ws.Range("O32").Value = 10.2
$ws.Range("P32").Value = 11.2
# Row 33: 全员劳动生产率
$ws.Range("K33").Value = "-"
$ws.Range("L33").Value = 41.8
$ws.Range("M33").Value = 40.5
$ws.Range("N33").Value = 40.2
$ws.Range("O33").Value = 44.6
$ws.Range("P33").Value = 50.7
# Row 34: 固定资产占营业总收入的比率
$ws.Range("K34").Value = 11.9
$ws.Range("L34").Value = 11
$ws.Range("M34").Value = 10.3
$ws.Range("N34").Value = 10.9
$ws.Range("O34").Value = 10.1
$ws.Range("P34").Value = 8.4
# Row 35: 职工薪酬占成本费用总额的比率
$ws.Range("K35").Value = 6.8
$ws.Range("L35").Value = 6.9
$ws.Range("M35").Value = 7.8
$ws.Range("N35").Value = 8.9
$ws.Range("O35").Value = 9.2
$ws.Range("P35").Value = 8.8

# --- Previously-blank placeholder cells within the existing 2018-2024 range ---
$ws.Range("E15").Value = " "
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "-"
$ws.Range("F25").Value = "-"
$ws.Range("G25").Value = "-"
$ws.Range("H25").Value = "-"

# --- View state: selection moved to L20 (and scrolled so column B leads) ---
try { $excel.ActiveWindow.ScrollColumn = 2 } catch { }
$null = $ws.Range("L20").Select()
